$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "备注"
$ws.Range("E2").Value = "空字符串"
$ws.Range("E3").Value = "备注1"
$ws.Range("E4").Value = "备注2"
$ws.Range("H4").Value = "是非得失"
$ws.Range("E5").Value = "备注3"
$ws.Range("E6").Value = "备注4"
$ws.Range("E7").Value = "备注5"
$ws.Range("E8").Value = "备注6"

$ws.Range("K8").Select() | Out-Null
